$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.299.37'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').Value = '1.899.77'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  -1.26%  '
$ws.Range('D5').Value = '''315.30'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').Value = '''1.004'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('D7').Value = '''0.5147'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '''0.3924'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('D9').Value = '''0.08452'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').Value = '''42.52'
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').Value = '''6.247'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').Value = '1.896.95'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '''20.77'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '''93.29'
$ws.Range('E17').Value = '  +1.87%  '
$ws.Range('D18').Value = '''0.00001107'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '''17.86'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').Value = '''6.029'
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').Value = '29.280.09'
$ws.Range('E23').Value = '  +2.21%  '
$ws.Range('D24').Value = '''11.17'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').Value = '''2.215'
$ws.Range('E25').Value = '  -2.61%  '
$ws.Range('D26').Value = '2.112.91'
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''159.06'
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''20.93'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').Value = '''2.443'
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('D30').Value = '''128.37'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').Value = '''1.060'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').Value = '''0.1047'
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('D33').Value = '''6.153'
$ws.Range('E33').Value = '  +5.90%  '
$ws.Range('D34').Value = '''3.661'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').Value = '''0.02474'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('D36').Value = '''0.06566'
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').Value = '''9.067'
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = '''1.234'
$ws.Range('E42').Value = '  -2.84%  '
$ws.Range('D43').Value = '''11.29'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '''0.6062'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').Value = '''13.19'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').Value = '''3.676'
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').Value = '''2.049'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('D48').Value = '''1.230'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').Value = '''123.53'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('E50').Value = '  -2.58%  '
$ws.Range('D51').Value = '''77.74'
$ws.Range('E51').Value = '  +0.59%  '
